$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "36.587.30"
Set-TextValue "E2" "  +2.28%  "
Set-TextValue "D3" "2.087.94"
Set-TextValue "E3" "  +10.47%  "
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "247.49"
Set-TextValue "E5" "  +0.53%  "
Set-TextValue "D6" "0.666"
Set-TextValue "E6" "  -3.78%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "D8" "45.08"
Set-TextValue "E8" "  +4.76%  "
Set-TextValue "D9" "60.62"
Set-TextValue "E9" "  +7.21%  "
Set-TextValue "E10" "  +1.85%  "
Set-TextValue "E11" "  -4.32%  "
Set-TextValue "E12" "  +0.41%  "
Set-TextValue "D13" "14.52"
Set-TextValue "E13" "  -2.62%  "
Set-TextValue "D14" "2.385.99"
Set-TextValue "E14" "  +10.15%  "
Set-TextValue "E15" "  +4.77%  "
Set-TextValue "D16" "2.074.93"
Set-TextValue "E16" "  +9.60%  "
Set-TextValue "D17" "4.92"
Set-TextValue "E17" "  -2.69%  "
Set-TextValue "D18" "36.611.66"
Set-TextValue "E18" "  +2.41%  "
Set-TextValue "D19" "71.79"
Set-TextValue "E19" "  -2.41%  "
Set-TextValue "D20" "0.0₃0815"
Set-TextValue "E20" "  -2.17%  "
Set-TextValue "D21" "238.14"
Set-TextValue "E21" "  -3.54%  "
Set-TextValue "D22" "12.74"
Set-TextValue "E22" "  -2.77%  "
Set-TextValue "D23" "4.93"
Set-TextValue "E23" "  -4.82%  "
Set-TextValue "E24" "  +0.09%  "
Set-TextValue "D25" "2.47"
Set-TextValue "E25" "  -8.93%  "
Set-TextValue "D26" "169.33"
Set-TextValue "E26" "  +1.95%  "
Set-TextValue "D27" "20.46"
Set-TextValue "E27" "  +10.84%  "
Set-TextValue "D28" "8.85"
Set-TextValue "E28" "  +2.09%  "
Set-TextValue "D29" "1.96"
Set-TextValue "E29" "  -8.81%  "
Set-TextValue "D30" "0.122"
Set-TextValue "E30" "  -5.12%  "
Set-TextValue "D31" "22.35"
Set-TextValue "E31" "  +58.68%  "
Set-TextValue "D32" "4.38"
Set-TextValue "E32" "  -1.17%  "
Set-TextValue "D33" "0.0584"
Set-TextValue "E33" "  -4.01%  "
Set-TextValue "D34" "0.0905"
Set-TextValue "E34" "  +17.27%  "
Set-TextValue "D35" "1.89"
Set-TextValue "E35" "  +1.34%  "
Set-TextValue "E36" "  -0.19%  "
Set-TextValue "D37" "2.30"
Set-TextValue "E37" "  +18.46%  "
Set-TextValue "E38" "  +4.83%  "
Set-TextValue "D39" "4.01"
Set-TextValue "E39" "  -6.32%  "
Set-TextValue "E40" "  -9.23%  "
Set-TextValue "D41" "1.15"
Set-TextValue "E41" "  +5.61%  "
Set-TextValue "D42" "97.77"
Set-TextValue "E42" "  -1.58%  "
Set-TextValue "D43" "0.0216"
Set-TextValue "E43" "  -5.93%  "
Set-TextValue "D44" "2.79"
Set-TextValue "E44" "  +15.81%  "
Set-TextValue "D45" "15.88"
Set-TextValue "E45" "  -6.38%  "
Set-TextValue "D46" "1.336.32"
Set-TextValue "E46" "  +1.46%  "
Set-TextValue "E47" "  +1.41%  "
Set-TextValue "E48" "  +3.85%  "
Set-TextValue "D49" "2.275.70"
Set-TextValue "E49" "  +10.05%  "
Set-TextValue "E50" "  -5.44%  "
Set-TextValue "D51" "3.87"
Set-TextValue "E51" "  +15.70%  "
